# Add a new "2022-Q3" sheet (fund holdings detail) positioned right after
# "总计", pushing "2022-Q1" and "2021-Q3" one slot later, and add the
# matching summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Force a value to be stored as TEXT (many of these columns hold
# numeric-looking strings like "9.01" or zero-padded codes like "002300"
# that must NOT be reinterpreted as numbers).
function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

function Set-NumCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" sheet by duplicating the "2022-Q1" sheet
#    (same column layout/styles) and placing the copy right after "总计".
# ---------------------------------------------------------------------

$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($null, $wb.Worksheets.Item(1))
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The template only has 4 rows (1 header + 3 data); we need 1 header + 5
# data rows, so extend the formatting of the last data row down two more
# rows before filling in values.
$newSheet.Range("A4:H4").Copy()
$newSheet.Range("A5:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$fundRows = @(
    @{A=0; B="002300"; C="长盛医疗行业量化配置股票";      D="2.39"; E="93.48"; F="5.79"; G="0.1384"; H=4},
    @{A=1; B="000684"; C="长盛养老健康产业灵活配置混合";    D="1.35"; E="92.15"; F="5.65"; G="0.0763"; H=3},
    @{A=2; B="008412"; C="长盛竞争优势股票A";              D="0.72"; E="91.35"; F="5.62"; G="0.0405"; H=3},
    @{A=3; B="008413"; C="长盛竞争优势股票C";              D="0.39"; E="91.35"; F="5.62"; G="0.0219"; H=3},
    @{A=4; B="006603"; C="嘉实互融精选股票";                D="0.12"; E="82.85"; F="2.87"; G="0.0034"; H=7}
)

$r = 2
foreach ($row in $fundRows) {
    Set-NumCell  $newSheet $r 1 $row.A
    Set-TextCell $newSheet $r 2 $row.B
    Set-TextCell $newSheet $r 3 $row.C
    Set-TextCell $newSheet $r 4 $row.D
    Set-TextCell $newSheet $r 5 $row.E
    Set-TextCell $newSheet $r 6 $row.F
    Set-TextCell $newSheet $r 7 $row.G
    Set-NumCell  $newSheet $r 8 $row.H
    $r++
}

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: the new quarter is inserted at the
#    top of the data (row 2), the old rows shift down one, and a fresh
#    row is appended for what used to be the last entry.
# ---------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

# Extend formatting from row 3 down to the new row 4.
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.28

$total.Cells.Item(3, 2).Value = "2022-Q1"
$total.Cells.Item(3, 3).Value = 3
$total.Cells.Item(3, 4).Value = 0.15

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q3"
$total.Cells.Item(4, 3).Value = 3
$total.Cells.Item(4, 4).Value = 0.6

# ---------------------------------------------------------------------
# 3) Keep "2021-Q3" as the active/selected tab, matching the original
#    workbook (the new sheet insertion would otherwise steal focus).
# ---------------------------------------------------------------------

$wb.Worksheets.Item("2021-Q3").Activate()
